# Update with some bolometers
# Adds new Mech (Standoff) and Connectors (DDR4 SODIMM, NGFF M.2, DF12NB, Molex, Amphenol) parts.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Mech sheet: fill in previously-blank row 11 with a new Standoff part
# ---------------------------------------------------------------------------
$mech = $wb.Worksheets.Item("Mech")

$mech.Range("A11").Value = "9774040360R"
$mech.Range("B11").Value = "9774040360R"
$mech.Range("C11").Value = "Standoff"
# A1:D1-style placeholder cells already exist (empty, unstyled) on this row;
# clear first so the write re-resolves the column's normal style.
$mech.Range("D11").ClearContents()
$mech.Range("D11").Value = "Common.SchLib"
$mech.Range("E11").Value = "ROUND STANDOFF M3X0.5 STEEL 4MM"
$mech.Range("F11").Value = "9774040360"
$mech.Range("H11").Value = "Main.PcbLib"
$mech.Range("I11").Value = "Würth Elektronik"
$mech.Range("J11").Value = "9774040360R"
$mech.Range("K11").Value = "Digi-Key"
$mech.Range("L11").Value = "732-5271-1-ND"

# ---------------------------------------------------------------------------
# Connectors sheet: append 7 new connector rows (47-53)
# ---------------------------------------------------------------------------
$conn = $wb.Worksheets.Item("Connectors")

# Row 47 - DDR4 SODIMM 260P socket
$conn.Range("A47").Value = "2309409-2"
$conn.Range("B47").Value = "2309409-2"
$conn.Range("C47").Value = "Connector P260 MH"
$conn.Range("D47").Value = "Connectors.SchLib"
$conn.Range("E47").Value = "DDR4 SODIMM 260P 5.2H STD"
$conn.Range("F47").Value = "JEDEC_DDR4_SODIMM_SKT_STD_H5.2"
$conn.Range("H47").Value = "Connectors.PcbLib"
$conn.Range("I47").Value = "TE Connectivity AMP Connectors"
$conn.Range("J47").Value = "2309409-2"
$conn.Range("K47").Value = "Digi-Key"
$conn.Range("L47").Value = "A141326CT-ND"

# Row 48 - NGFF M.2 Key E connector
$conn.Range("A48").Value = "1-2199230-1"
$conn.Range("B48").Value = "Connector NGFF M.2 Key E"
$conn.Range("C48").Value = "Connector NGFF M.2 Key E"
$conn.Range("D48").Value = "Connectors.SchLib"
$conn.Range("E48").Value = "PCI Express/PCI Connectors M.2 0.5PITCH 4.2H KEY E 10U AU"
$conn.Range("F48").Value = "TE_1-2199230-1"
$conn.Range("H48").Value = "Connectors.PcbLib"
$conn.Range("I48").Value = "TE Connectivity AMP Connectors"
$conn.Range("J48").Value = "1-2199230-1"
$conn.Range("K48").Value = "Mouser"
$conn.Range("L48").Value = "571-1-2199230-1"

# Row 49 - NGFF M.2 Key B connector
$conn.Range("A49").Value = "1-2199230-0"
$conn.Range("B49").Value = "Connector NGFF M.2 Key B"
$conn.Range("C49").Value = "Connector NGFF M.2 Key B"
$conn.Range("D49").Value = "Connectors.SchLib"
$conn.Range("E49").Value = "PCI Express/PCI Connectors M.2 0.5PITCH 4.2H KEY B 10U AU"
$conn.Range("F49").Value = "TE_1-2199230-0"
$conn.Range("H49").Value = "Connectors.PcbLib"
$conn.Range("I49").Value = "TE Connectivity AMP Connectors"
$conn.Range("J49").Value = "1-2199230-0"
$conn.Range("K49").Value = "Mouser"
$conn.Range("L49").Value = "571-1-2199230-0"

# Row 50 - Hirose DF12NB(4.0) 50-pos header
$conn.Range("A50").Value = "DF12NB(4.0)-50DP-0.5V(51)"
$conn.Range("B50").Value = "DF12NB(4.0)-50DP-0.5V(51)"
$conn.Range("C50").Value = "Connector P50"
$conn.Range("D50").Value = "Connectors.SchLib"
$conn.Range("E50").Value = "CONN HDR 50POS SMD GOLD"
$conn.Range("F50").Value = "HIROSE_DF12NB(4.0)-50DP-0.5V"
$conn.Range("H50").Value = "Connectors.PcbLib"
$conn.Range("I50").Value = "Hirose Electric Co Ltd"
$conn.Range("J50").Value = "DF12NB(4.0)-50DP-0.5V(51)"
$conn.Range("K50").Value = "Digi-Key"
$conn.Range("L50").Value = "26-DF12NB(4.0)-50DP-0.5V(51)CT-ND"

# Row 51 - Hirose DF12NB 50-pos receptacle
$conn.Range("A51").Value = "DF12NB-50DS-0.5V(51)"
$conn.Range("B51").Value = "DF12NB-50DS-0.5V(51)"
$conn.Range("C51").Value = "Connector P50"
$conn.Range("D51").Value = "Connectors.SchLib"
$conn.Range("E51").Value = "CONN RCPT 50POS SMD GOLD"
$conn.Range("F51").Value = "HIROSE_DF12NB-50DS-0.5V"
$conn.Range("H51").Value = "Connectors.PcbLib"
$conn.Range("I51").Value = "Hirose Electric Co Ltd"
$conn.Range("J51").Value = "DF12NB-50DS-0.5V(51)"
$conn.Range("K51").Value = "Digi-Key"
$conn.Range("L51").Value = "26-DF12NB-50DS-0.5V(51)CT-ND"

# Row 52 - Molex 10-pos header
$conn.Range("A52").Value = "0430451027"
$conn.Range("B52").Value = "0430451027"
$conn.Range("C52").Value = "Connector P10"
$conn.Range("D52").Value = "Connectors.SchLib"
$conn.Range("E52").Value = "CONN HEADER VERT 10POS 3MM"
$conn.Range("F52").Value = "MOLEX_430451027"
$conn.Range("H52").Value = "Connectors.PcbLib"
$conn.Range("I52").Value = "Molex"
$conn.Range("J52").Value = "0430451027"
$conn.Range("K52").Value = "Digi-Key"
$conn.Range("L52").Value = "WM7488-ND"

# Row 53 - Amphenol RJ45 jack
$conn.Range("A53").Value = "94152-088LF"
$conn.Range("B53").Value = "94152-088LF"
$conn.Range("C53").Value = "Connector P8"
$conn.Range("D53").Value = "Connectors.SchLib"
$conn.Range("E53").Value = "CONN MOD JACK 8P8C VERT SHIELDED"
$conn.Range("F53").Value = "AMPHENOL_94152-088LF"
$conn.Range("H53").Value = "Connectors.PcbLib"
$conn.Range("I53").Value = "Amphenol ICC (FCI)"
$conn.Range("J53").Value = "94152-088LF"
$conn.Range("K53").Value = "Digi-Key"
$conn.Range("L53").Value = "609-1072-ND"
